$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "A2" = "华胜天成"; "B2" = "华胜天成"; "C2" = "华夏幸福";
    "A3" = "蓝色光标"; "B3" = "华夏幸福"; "C3" = "航天发展";
    "A4" = "孚日股份"; "B4" = "孚日股份"; "C4" = "华胜天成";
    "A5" = "华夏幸福"; "B5" = "三六零"; "C5" = "孚日股份";
    "A6" = "浪潮软件"; "B6" = "航天发展"; "C6" = "海马汽车";
    "A7" = "平潭发展"; "B7" = "蓝色光标"; "C7" = "多氟多";
    "A8" = "航天发展"; "B8" = "利欧股份"; "C8" = "人民同泰";
    "A9" = "福石控股"; "B9" = "雪人集团"; "C9" = "利欧股份";
    "A10" = "利欧股份"; "B10" = "浪潮软件"; "C10" = "雪人集团";
    "B11" = "长城军工"; "C11" = "平潭发展";
    "A12" = "海马汽车"; "B12" = "多氟多"; "C12" = "安泰集团";
    "B13" = "海马汽车"; "C13" = "浪潮软件";
    "A14" = "多氟多"; "B14" = "福石控股"; "C14" = "蓝色光标";
    "A15" = "人民同泰"; "B15" = "平潭发展"; "C15" = "首开股份";
    "A16" = "安泰集团"; "B16" = "首开股份"; "C16" = "胜利股份";
    "A17" = "榕基软件"; "B17" = "人民同泰"; "C17" = "龙洲股份";
    "A18" = "首开股份"; "B18" = "视觉中国"; "C18" = "三六零";
    "A19" = "宣亚国际"; "B19" = "安泰集团"; "C19" = "常山北明";
    "A20" = "视觉中国"; "B20" = "海南海药"; "C20" = "九牧王";
    "A21" = "龙洲股份"; "B21" = "日出东方"; "C21" = "天际股份";
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
